$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 32

$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 9
$ws.Range("C7").Value = 12
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 18
$ws.Range("C10").Value = 21
$ws.Range("C11").Value = 24
$ws.Range("C12").Value = 27
$ws.Range("C13").Value = 30

$excel.Calculate()
